$wb = $excel.ActiveWorkbook

# --- Add the new "When_Dup" worksheet as the LAST sheet in the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "When_Dup"

# --- Header row ---
$ws.Range("A1").Value = "Item_Grade"
$ws.Range("B1").Value = "Acquired_Mileage"

# --- Data rows ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 10
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 20
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 30

# Copy the header style used by the other sheets (bold font on yellow fill,
# centered) onto the new header cells.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0 | Out-Null

# --- Column widths (engine persists ColumnWidth + 5/7 as the stored <col> width) ---
$ws.Columns("A").ColumnWidth = 12 - 5/7
$ws.Columns("B").ColumnWidth = 18 - 5/7

# --- Threaded comment on A1 describing the grade values ---
$ws.Range("A1").AddCommentThreaded("Normal = 1`nRare = 2`nUnique = 3") | Out-Null

# --- Selection state matching the authored file ---
$ws.Range("F21").Select() | Out-Null

Write-Output "done"
